# Apply the commit's changes to the workbook:
#  - rename the sheet "Under Grad" -> "Comp_student Numbers_Oct 2022"
#  - narrow column F, merge/narrow columns G:H to a uniform width, and narrow column K
#
# Note: the engine stores ColumnWidth on an integer-pixel (1/6 character) grid,
# so the input "chars" value is chosen so the stored width lands on (or as close
# as possible to) the target raw OOXML width after the engine's +5/6 padding.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet/tab
$ws.Name = "Comp_student Numbers_Oct 2022"

# Column F (6): width 26.5703125 -> 20
$ws.Columns.Item(6).ColumnWidth = 19.166666666666668

# Columns G:H (7:8): widths 11.28515625 / 10 -> both become 10
$ws.Range("G1:H1").EntireColumn.ColumnWidth = 9.166666666666666

# Column K (11): width 40.28515625 -> 36.1640625
$ws.Columns.Item(11).ColumnWidth = 35.330729166666664
